$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bad Drivers table (rows 3-7) ---
$ws.Range("A3").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.220.1.1"
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 2424
$ws.Range("D3").Value = 93.8

$ws.Range("A4").Value = "Intel(R) Dual Band Wireless-AC 7265 - 19.51.48.1"
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 478
$ws.Range("D4").Value = 98.09999999999999

$ws.Range("A5").Value = "Intel(R) Dual Band Wireless-AC 7260 - 17.15.0.5"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 296
$ws.Range("D5").Value = 98.3

$ws.Range("A6").Value = "Intel(R) Dual Band Wireless-AC 7265 - 19.51.42.2"
$ws.Range("B6").Value = 39
$ws.Range("C6").Value = 3175
$ws.Range("D6").Value = 98.5

$ws.Range("A7").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.32.1"
$ws.Range("B7").Value = 8
$ws.Range("C7").Value = 703
$ws.Range("D7").Value = 98.5

# --- Totals row (row 8) ---
$ws.Range("C8").Value = 7076

# --- Good Drivers table (rows 16-36) ---
# Keep column E as text (it holds literal date-strings, not real dates)
$ws.Range("E16:E36").NumberFormat = "@"

$ws.Range("A16").Value = "Intel(R) Dual Band Wireless-AC 7265 - 19.50.1.6"
$ws.Range("B16").Value = 29259
$ws.Range("D16").Value = 99.90000000000001

$ws.Range("A17").Value = "Intel(R) Dual Band Wireless-AC 7265 - 19.51.8.3"
$ws.Range("B17").Value = 13028
$ws.Range("D17").Value = 100

$ws.Range("A18").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B18").Value = 445055
$ws.Range("D18").Value = 99.90000000000001
$ws.Range("E18").Value = "2024-11-10"

$ws.Range("A19").Value = "Intel(R) Dual Band Wireless-AC 7265 - 19.51.50.2"
$ws.Range("B19").Value = 19910
$ws.Range("D19").Value = 100
$ws.Range("E19").Value = "2023-11-06"

$ws.Range("A20").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B20").Value = 77849
$ws.Range("D20").Value = 99.90000000000001
$ws.Range("E20").Value = "2021-08-18"

$ws.Range("A21").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B21").Value = 34244
$ws.Range("D21").Value = 100
$ws.Range("E21").Value = "2021-04-27"

$ws.Range("A22").Value = "Intel(R) Dual Band Wireless-AC 8265 - 22.30.0.11"
$ws.Range("B22").Value = 170510
$ws.Range("D22").Value = 99.90000000000001
$ws.Range("E22").Value = "2021-01-19"

$ws.Range("A23").Value = "Intel(R) Dual Band Wireless-AC 8265 - 22.0.1.1"
$ws.Range("B23").Value = 52096
$ws.Range("D23").Value = 100
$ws.Range("E23").Value = "2020-09-28"

$ws.Range("A24").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B24").Value = 59673
$ws.Range("D24").Value = 100
$ws.Range("E24").Value = "2020-08-05"

$ws.Range("A25").Value = "Intel(R) Dual Band Wireless-AC 7265 - 19.51.30.1"
$ws.Range("B25").Value = 201061
$ws.Range("D25").Value = 100
$ws.Range("E25").Value = "2020-06-01"

$ws.Range("A26").Value = "Intel(R) Dual Band Wireless-AC 7265 - 19.51.29.1"
$ws.Range("B26").Value = 40159
$ws.Range("D26").Value = 100
$ws.Range("E26").Value = "2020-04-15"

$ws.Range("A27").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B27").Value = 113652
$ws.Range("D27").Value = 100
$ws.Range("E27").Value = "2020-01-06"

$ws.Range("A28").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B28").Value = 56018
$ws.Range("D28").Value = 100
$ws.Range("E28").Value = "2019-12-14"

$ws.Range("A29").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.11.3"
$ws.Range("B29").Value = 161874
$ws.Range("D29").Value = 100
$ws.Range("E29").Value = "2019-09-05"

$ws.Range("A30").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.12.5"
$ws.Range("B30").Value = 143342
$ws.Range("D30").Value = 99.90000000000001
$ws.Range("E30").Value = "2019-08-25"

$ws.Range("A31").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.10.2"
$ws.Range("B31").Value = 20227
$ws.Range("D31").Value = 100
$ws.Range("E31").Value = "2019-05-11"

$ws.Range("A32").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.9.1"
$ws.Range("B32").Value = 34065
$ws.Range("D32").Value = 100
$ws.Range("E32").Value = "2019-04-28"

$ws.Range("A33").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.8.1"
$ws.Range("B33").Value = 48540
$ws.Range("D33").Value = 100
$ws.Range("E33").Value = "2019-03-16"

$ws.Range("A34").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.5.2"
$ws.Range("B34").Value = 184564
$ws.Range("D34").Value = 99.90000000000001
$ws.Range("E34").Value = "2018-11-25"

$ws.Range("A35").Value = "Intel(R) Dual Band Wireless-AC 7260 - 18.33.15.1"
$ws.Range("B35").Value = 83189
$ws.Range("D35").Value = 100
$ws.Range("E35").Value = "2018-11-10"

$ws.Range("A36").Value = "Intel(R) Dual Band Wireless-AC 7265 - 19.51.14.1"
$ws.Range("B36").Value = 120862
$ws.Range("D36").Value = 100
$ws.Range("E36").Value = "2018-05-26"
